# Update the "Marking" row (Right count) and the "Total" row (Right count
# and the Correct/Total marks summary) on the "quiz" marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row - Right column
$ws.Range("B11").Value = 5

# Total row - Right column, and Corr/total marks label
$ws.Range("B12").Value = 90
$ws.Range("E12").Value = "90/140"
